$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 91820377
$ws.Range("B2").Value = 77506
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 792384.9250397275
$ws.Range("R2").Value = 7354586.848309446

$ws.Range("A3").Value = 91820374
$ws.Range("B3").Value = 77506
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 792384.1932845772
$ws.Range("R3").Value = 7354578.223190884

$ws.Range("A4").Value = 91820376
$ws.Range("B4").Value = 77506
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 792389.2480045473
$ws.Range("R4").Value = 7354564.908565938

$ws.Range("A5").Value = 91820391
$ws.Range("B5").Value = 78503
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 6456
$ws.Range("F5").Value = "Skinnlav"
$ws.Range("G5").Value = "Leptogium saturninum"
$ws.Range("H5").Value = "(Dicks.) Nyl."
$ws.Range("Q5").Value = 792447.929428296
$ws.Range("R5").Value = 7354525.833729008

$ws.Range("A6").Value = 91820347
$ws.Range("B6").Value = 78569
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6458
$ws.Range("F6").Value = "Lunglav"
$ws.Range("G6").Value = "Lobaria pulmonaria"
$ws.Range("H6").Value = "(L.) Hoffm."
$ws.Range("Q6").Value = 792501.2096894301
$ws.Range("R6").Value = 7354565.988424786

$ws.Range("A7").Value = 91820325
$ws.Range("B7").Value = 78569
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6458
$ws.Range("F7").Value = "Lunglav"
$ws.Range("G7").Value = "Lobaria pulmonaria"
$ws.Range("H7").Value = "(L.) Hoffm."
$ws.Range("Q7").Value = 792453.9552827136
$ws.Range("R7").Value = 7354651.845769764

$ws.Range("A8").Value = 91820354
$ws.Range("B8").Value = 89410
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = "Granticka"
$ws.Range("G8").Value = "Porodaedalea chrysoloma"
$ws.Range("H8").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q8").Value = 792558.2461457669
$ws.Range("R8").Value = 7354663.934873462

$ws.Range("A9").Value = 91820321
$ws.Range("B9").Value = 78569
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 792498.060890534
$ws.Range("R9").Value = 7354686.158366068

$ws.Range("A10").Value = 91820352
$ws.Range("B10").Value = 78569
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6458
$ws.Range("F10").Value = "Lunglav"
$ws.Range("G10").Value = "Lobaria pulmonaria"
$ws.Range("H10").Value = "(L.) Hoffm."
$ws.Range("Q10").Value = 792451.0675179539
$ws.Range("R10").Value = 7354667.828167814

$ws.Range("A11").Value = 91820382
$ws.Range("B11").Value = 78503
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 6456
$ws.Range("F11").Value = "Skinnlav"
$ws.Range("G11").Value = "Leptogium saturninum"
$ws.Range("H11").Value = "(Dicks.) Nyl."
$ws.Range("Q11").Value = 792439.9662732746
$ws.Range("R11").Value = 7354520.118926036

$ws.Range("A12").Value = 91820367
$ws.Range("B12").Value = 78569
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6458
$ws.Range("F12").Value = "Lunglav"
$ws.Range("G12").Value = "Lobaria pulmonaria"
$ws.Range("H12").Value = "(L.) Hoffm."
$ws.Range("Q12").Value = 792462.0294626241
$ws.Range("R12").Value = 7354500.026836612

$ws.Range("A13").Value = 91820383
$ws.Range("B13").Value = 77506
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 792405.1790341048
$ws.Range("R13").Value = 7354544.991983407

$ws.Range("A14").Value = 91820369
$ws.Range("B14").Value = 78472
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 388
$ws.Range("F14").Value = "Stiftgelélav"
$ws.Range("G14").Value = "Collema furfuraceum"
$ws.Range("H14").Value = "(Arnold) Du Rietz"
$ws.Range("Q14").Value = 792412.1402423121
$ws.Range("R14").Value = 7354529.026110045

